$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "75.649.74"
$ws.Range("E2").Value = "  +8.44%  "
$ws.Range("D3").Value = "2.718.46"
$ws.Range("E3").Value = "  +11.61%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "187.04"
$ws.Range("E5").Value = "  +11.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "590.97"
$ws.Range("E6").Value = "  +4.38%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.542"
$ws.Range("E8").Value = "  +5.24%  "
$ws.Range("E9").Value = "  +14.17%  "
$ws.Range("D10").Value = "2.717.12"
$ws.Range("E10").Value = "  +11.62%  "
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.362"
$ws.Range("E12").Value = "  +8.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.78"
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").Value = "3.219.17"
$ws.Range("E14").Value = "  +11.65%  "
$ws.Range("D15").Value = "75.499.24"
$ws.Range("E15").Value = "  +8.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000190"
$ws.Range("E16").Value = "  +6.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.02"
$ws.Range("E17").Value = "  +12.18%  "
$ws.Range("D18").Value = "2.719.33"
$ws.Range("E18").Value = "  +11.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.35"
$ws.Range("E19").Value = "  +29.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.15"
$ws.Range("E20").Value = "  +11.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "378.42"
$ws.Range("E21").Value = "  +9.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.30"
$ws.Range("E22").Value = "  +14.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.11"
$ws.Range("E23").Value = "  +6.19%  "
$ws.Range("E24").Value = "  +4.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.13"
$ws.Range("E25").Value = "  +7.63%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.22"
$ws.Range("E27").Value = "  +10.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.61"
$ws.Range("E28").Value = "  +12.75%  "
$ws.Range("D29").Value = "2.860.07"
$ws.Range("E29").Value = "  +11.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("D31").Value = "0.0₃0990"
$ws.Range("E31").Value = "  +15.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "522.28"
$ws.Range("E32").Value = "  +14.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.41"
$ws.Range("E33").Value = "  +12.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.89"
$ws.Range("E34").Value = "  +6.72%  "
$ws.Range("E35").Value = "  +10.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  +7.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.14"
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.59"
$ws.Range("E39").Value = "  +7.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.37"
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "173.74"
$ws.Range("E42").Value = "  +27.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.07"
$ws.Range("E43").Value = "  +14.40%  "
$ws.Range("E44").Value = "  +12.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.334"
$ws.Range("E45").Value = "  +9.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.22"
$ws.Range("E46").Value = "  +12.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.40"
$ws.Range("E47").Value = "  +14.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "39.21"
$ws.Range("E48").Value = "  +2.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0858"
$ws.Range("E49").Value = "  +18.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.71"
$ws.Range("E50").Value = "  +9.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.549"
$ws.Range("E51").Value = "  +11.61%  "
